$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 145, shifting existing rows 145-257 down to 146-258
$ws.Rows.Item(145).Insert()

# Populate the newly inserted row 145 with the new record's data
$ws.Range("A145").Value = 3
$ws.Range("B145").Value = "Femacal de La Calera"
$ws.Range("C145").Value = "Coquimbo"
$ws.Range("D145").Value = 44978
$ws.Range("E145").Value = 5
$ws.Range("F145").Value = 100112030
$ws.Range("G145").Value = "Poroto granado"
$ws.Range("H145").Value = "Sin especificar"
$ws.Range("I145").Value = "Primera"
$ws.Range("J145").Value = 73
$ws.Range("K145").Value = 25000
$ws.Range("L145").Value = 26000
$ws.Range("M145").Value = 25479
$ws.Range("N145").Value = "$/malla 25 kilos"
$ws.Range("O145").Value = "Provincia de Limar" + [char]0x00ED
$ws.Range("P145").Value = 1019
$ws.Range("Q145").Value = 25
$ws.Range("R145").Value = "Hortaliza"
